$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Right hand table (A1:F7) -------------------------------------------
# Fingers that close towards 0 / open towards 180
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "close"
$ws.Range("E2").Value = 180
$ws.Range("F2").Value = "open"

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "close"
$ws.Range("E3").Value = 180
$ws.Range("F3").Value = "open"

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = "close"
$ws.Range("E4").Value = 180
$ws.Range("F4").Value = "open"

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "close"
$ws.Range("E5").Value = 180
$ws.Range("F5").Value = "open"

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = "open"
$ws.Range("E6").Value = 180
$ws.Range("F6").Value = "close"

# --- Left hand table (A9:F15) -------------------------------------------
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = "open"
$ws.Range("E10").Value = 180
$ws.Range("F10").Value = "close"

$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "open"
$ws.Range("E11").Value = 180
$ws.Range("F11").Value = "close"

$ws.Range("C12").Value = 0
$ws.Range("D12").Value = "open"
$ws.Range("E12").Value = 180
$ws.Range("F12").Value = "close"

$ws.Range("C13").Value = 0
$ws.Range("D13").Value = "open"
$ws.Range("E13").Value = 180
$ws.Range("F13").Value = "close"

$ws.Range("C14").Value = 0
$ws.Range("D14").Value = "open"
$ws.Range("E14").Value = 180
$ws.Range("F14").Value = "close"

# Header row additions: a highlighted note marker + note text
$ws.Range("H1").Value = "*Note"
$ws.Range("H1").Interior.Color = 65535
$ws.Range("I1").Value = "Apply roughly ~7 V"

[void]$ws.Range("L6").Select()
